# Reduced to 3 stim, added gray circles
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the 2nd rating row (was "ratingCS-1.png") to "ratingCS+3.png"
$ws.Range("A3").Value = "Instructions_EN/ratingCS+3.png"

# 2. Add a new 4th stim row "ratingCS+4.png" / -0.1, copying row 2's look
$ws.Range("A4").Value = "Instructions_EN/ratingCS+4.png"
$ws.Range("B4").Value = -0.1

# 3. Give row 2 (the alternating "stim" row) a light gray fill (White,
#    Background 1, Darker 5% -> RGB 242,242,242), and make the new row 4
#    match it.
$ws.Range("A2:B2").Interior.Color = 15921906

$ws.Range("A4:B4").Interior.Color = 15921906

# 4. Row heights: row2 -> 15 (customHeight), row4 -> 16 (matches original row height)
$ws.Rows.Item(2).RowHeight = 15
$ws.Rows.Item(4).RowHeight = 16

$ws.Range("A11").Select()
